# Adds a new "Kangaroo mother care" program as row 16 (alphabetically
# between "IYCF 1" and "Lipid-based nutrition supplements") to the three
# per-program sheets: "Programs to include", "Coverage scenario" and
# "Budget scenario". Existing rows 16-32 shift down to 17-33.

$wb = $excel.ActiveWorkbook

function Add-KangarooRow($sheetName, $lastCol, $bLabel) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new blank row at 16 - shifts rows 16:32 down to 17:33,
    # along with any values/formatting they carried (sorted list stays
    # intact apart from the new entry).
    $ws.Rows.Item(16).Insert()

    # Copy the (now shifted-down) row 17 formatting onto the blank row 16
    # so the new row matches its neighbours exactly (only columns A:lastCol,
    # to avoid touching the whole 16384-column row).
    $srcRange = "A17:" + $lastCol + "17"
    $dstRange = "A16:" + $lastCol + "16"
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)

    # Fill in the new program row.
    $ws.Range("A16").Value = "Kangaroo mother care"
    $ws.Range("B16").Value = $bLabel

    # Refresh the sheet's recorded sort state (A2:..32 -> A2:..33) the same
    # way Excel does when you re-run Data > Sort over the grown range.
    $sortRange = "A2:" + $lastCol + "33"
    $keyRange = "A2:A33"
    $sortObj = $ws.Sort
    $sortObj.SortFields.Clear()
    $sortObj.SortFields.Add($ws.Range($keyRange))
    $sortObj.SetRange($ws.Range($sortRange))
    $sortObj.Header = -4142
    $sortObj.Apply()

    # Leave the selection on the newly inserted cell, matching the final
    # on-screen state after the edit.
    $ws.Range("B17").Select()
}

Add-KangarooRow "Programs to include" "B" "x"
Add-KangarooRow "Coverage scenario" "K" "Coverage"
Add-KangarooRow "Budget scenario" "K" "Spending"

# "Programs to include" ends up as the active tab.
$wb.Worksheets.Item("Programs to include").Activate()
$wb.Worksheets.Item("Programs to include").Range("B17").Select()
